$wb = $excel.ActiveWorkbook

# Sheet: VENTAS POR GRUPO
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("M12").Value = 704.5
$wsVentasGrupo.Range("M20").Value = 2680.32

# Sheet: VENTA MENSUAL
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F12").Value = 894.58
$wsVentaMensual.Range("F20").Value = 3135.75
$wsVentaMensual.Range("F36").Value = 13812.23

# Sheet: CUMPLIMIENTO MENSUAL
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D12").Value = 12501.42
$wsCumplimiento.Range("E12").Value = 9199.85
$wsCumplimiento.Range("F12").Value = 0.5760685895341608
$wsCumplimiento.Range("D14").Value = 13812.23
$wsCumplimiento.Range("E14").Value = 22773.33723718183
$wsCumplimiento.Range("F14").Value = 0.3775322085470541
